$d = $word.ActiveDocument

# --- Paragraph 1 (title): two runs separated by a line break; title text itself has no quotes, so Find/Replace is safe here ---
$d.Content.Find.Execute('המאמר היומי של מייק - 10.02.25', $true, $false, $false, $false, $false, $true, 1, $false, 'המאמר היומי של מייק - 09.02.25', 2) | Out-Null
$d.Content.Find.Execute('On the expressiveness and spectral bias of KANs', $true, $false, $false, $false, $false, $true, 1, $false, 'Why Is Anything Conscious?', 2) | Out-Null

# --- Paragraphs 3-15: direct text replacement via Range.Text (avoids smart-quote autocorrect on straight quotes/apostrophes) ---
$d.Paragraphs(3).Range.Text = 'המאמר המעניין מאת מייקל טימותי בנט, שון וולש ואנה צ''יאוניקה מתמודד עם "הבעיה הקשה של התודעה", שנוסחה על ידי דייויד צ''אלמרס(David John Chalmers). אתגר פילוסופי זה מעלה את השאלה מדוע עיבוד מידע במערכות מסוימות, במיוחד ביולוגיות, מוביל לחוויות סובייקטיביות או *קוואליה*. המחברים מציעים שינוי פרדיגמה, המעגן את התודעה בדינמיקה של מערכות self-organizing שעוצבו על ידי הברירה הטבעית.'
$d.Paragraphs(4).Range.Text = 'הם טוענים כי תודעה תופעתית (phenomenal) - החוויה הסובייקטיבית של "איך זה מרגיש" - אינה רק יסודית אלא הכרחית להתנהגות אדפטיבית. מעניין כי באמצעות פריימוורק חישובי פורמלי, המחברים טוענים נגד האפשרות של "זומבים", מערכות המתפקדות כמו בני אדם אך חסרות חוויה סובייקטיבית, ומצהירים באופן פרובוקטיבי כי "הטבע אינו אוהב זומבים". חוויה סובייקטיבית היא ההבנה המלאה והחווייתית של ההשפעה הרגשית והקוגניטיבית כאחד הנובעת מאופן שבו הבני אדם מבינים ומפרשים אירועים שנצפו או נחוו על ידי הם.'
$d.Paragraphs(5).Range.Text = 'תרומות מרכזיות: '
$d.Paragraphs(6).Range.Text = 'מסגרת מתמטית לאנקטיביזם פנ-חישובי'
$d.Paragraphs(7).Range.Text = 'המחברים מציגים מערכת פורמלית המעוגנת ב*פנ-חישוביות* ו*אנקטיביזם*(Pancomputational Enactivism). פנ-חישוביות מניחה שכל המערכות הדינמיות מחשבות משהו, בעוד שאנקטיביזם מדגיש את ההכרה כנובעת מאינטראקציות בין מערכת לסביבתה. האלמנטים המרכזיים במודל שלהם כוללים:'
$d.Paragraphs(8).Range.Text = '- סביבה: מוגדרת כקבוצת מצבים, עם מעברים המתוארים על ידי תכנות דקלרטיבי.'
$d.Paragraphs(9).Range.Text = '- שכבת הפשטה: מבנה המגדיר כיצד מערכות מפרשות היבטים סביבתיים.'
$d.Paragraphs(10).Range.Delete()
$d.Paragraphs(9).Range.InsertParagraphAfter()
$d.Paragraphs(10).Range.Text = '- משימות ומדיניות: מבני התנהגות הממפים קלט לפלט, המאפשרים התנהגות אדפטיבית.'
$d.Paragraphs(11).Range.Text = '- זהויות סיבתיות:*ייצוגים של התערבויות והשפעותיהן, חיוניים למודעות עצמית.'
$d.Paragraphs(12).Range.Text = 'הפריימוורק מתאר כיצד מערכות מודעות שומרות על קוהרנטיות והסתגלות על ידי בניית זהויות סיבתיות מורכבות יותר ויותר, המהוות בסיס למודעות עצמית.'
$d.Paragraphs(13).Range.Text = 'היררכיה של תודעה'
$d.Paragraphs(14).Range.Text = 'תובנה מרכזית היא ההתפתחות ההיררכית של התודעה, המונעת על ידי ברירה טבעית ולחצי סקאלה. המחברים מתארים 6 שלבים מתקדמים:'
$d.Paragraphs(15).Range.Text = '1. מערכות לא מודעות: ישויות חסרות חוויה או הכרה, כמו סלעים.'

# --- Insert 17 new paragraphs after paragraph 15 (now holding the "1. מערכות לא מודעות..." text) ---
# Insert empty paragraphs one at a time, always right after paragraph 15, so they end up in forward order.
for ($i = 0; $i -lt 17; $i++) {
    $d.Paragraphs(15 + $i).Range.InsertParagraphAfter()
}

# Fill in the text for the newly created paragraphs 16..32
$d.Paragraphs(16).Range.Text = '2. מערכות מקודדות באופן קשיח: מערכות עם תגובות קבועות, מתוכנתות מראש (למשל, חד-תאיים).'
$d.Paragraphs(17).Range.Text = '3. מערכות לומדות: מערכות מסתגלות ללא מודעות עצמית (למשל, תולעים נמטודות).'
$d.Paragraphs(18).Range.Text = '4. מערכות עצמי מסדר ראשון: מסוגלות להבחין בין פעולות שנוצרו עצמאית לבין אירועים חיצוניים (למשל, זבובי בית).'
$d.Paragraphs(19).Range.Text = '5. מערכות עצמיות מסדר שני: מסוגלות למטא-ייצוג ותקשורת מכוונת (למשל, עורבים).'
$d.Paragraphs(20).Range.Text = '6. מערכות עצמי מסדר שלישי: ישויות רפלקטיביות במלואן המסוגלות לחשוב על המודעות שלהן עצמן (למשל, בני אדם).'
$d.Paragraphs(21).Range.Text = 'היררכיה זו מדגישה כיצד היבטים איכותיים של תודעה מתפתחים באופן טבעי ככל שמערכות נעשות מסוגלות יותר למדל את עצמן ואת סביבתן.'
$d.Paragraphs(22).Range.Text = 'עיבוד איכותי וכמותי:'
$d.Paragraphs(23).Range.Text = 'המחברים טוענים כי *איכות קודמת לכמות* בעיבוד מידע. לפני שאורגניזם יכול לתייג או למדוד מידע, עליו לחוות הבדלים איכותיים. תודעה פנומנלית מתפתחת מכיוון שמערכות חיות חייבות לסווג ולתעדף מידע הרלוונטי להישרדות. סיווגים איכותיים אלה מהווים את הבסיס לחוויה סובייקטיבית. טענה זו מאתגרת תיאוריות חישוביות מסורתיות, המתייחסות לעתים קרובות לתודעה כתהליך ייצוגי טהור. על ידי הדגשת הקדימות של החוויה האיכותית, המחברים מספקים פרספקטיבה רעננה על מקורות התודעה.'
$d.Paragraphs(24).Range.Text = 'גישת עקרונות ראשוניים:'
$d.Paragraphs(25).Range.Text = 'הפורמליזם במאמר נגזר משתי אקסיומות בסיסיות:'
$d.Paragraphs(26).Range.Text = '1. במקום שיש דברים, אנו קוראים לדברים אלה הסביבה.'
$d.Paragraphs(27).Range.Text = '2. במקום שדברים שונים, יש לנו מצבים שונים של הסביבה.'
$d.Paragraphs(28).Range.Text = 'אקסיומות אלה מובילות לצורה חסרת ייצוג של פנ-חישוביות, בה מצבים ומעברים מגדירים סביבות מבלי להניח מבנים פנימיים ספציפיים. המחברים ממסגרים ארגון עצמי כיכולת להגביל פלטים על בסיס קלטים, ובכך להשיג התנהגות אדפטיבית.'
$d.Paragraphs(29).Range.Text = 'דחיית זומבים'
$d.Paragraphs(30).Range.Text = 'אחת הטענות המעניינות ביותר במאמר היא ש"הטבע אינו אוהב זומבים". המחברים טוענים שתודעה פנומנלית חיונית למודעות גישה ולהתנהגות אדפטיבית. תוכן ייצוגי - מה שאורגניזמים חושבים עליו - נגזר תמיד מחוויה איכותית. לכן, מערכת המתנהגת כמו ישות מודעת חייבת בהכרח לחוות חוויה סובייקטיבית. טענה זו מאתגרת ישירות ניסויי מחשבה המציעים את קיומן של ישויות לא מודעות אך זהות בהתנהגותן.'
$d.Paragraphs(31).Range.Text = 'קשרים אמפיריים'
$d.Paragraphs(32).Range.Text = 'המאמר מבוסס על ממצאים אמפיריים לגבי *רה-אפרנציה*, כלומר היכולת להבחין בין גירויים שנוצרו עצמאית לבין גירויים חיצוניים. רה-אפרנציה, הנצפית ביונקים וחרקים, קשורה ליצירת עצמי מסדר ראשון. המחברים גוזרים מבנה זה מעקרונות מתמטיים ומיישרים את מסקנותיהם עם עבודתם של מרקר, ברון וקליין.'

# --- Remaining paragraphs: סיכום: (Word para 33) stays unchanged; summary body + URL replaced ---
$d.Paragraphs(34).Range.Text = 'המאמר מציע גישה מסקרנת לבעיה הקשה של התודעה על ידי עיגונה בברירה טבעית, ארגון עצמי ופורמליזם חישובי. המסגרת ההיררכית של המחברים מספקת הסבר משכנע לאופן שבו תודעה מתפתחת ומדוע חוויה סובייקטיבית היא יסודית להתנהגות אדפטיבית. טענתם הפרובוקטיבית שזומבים הם בלתי אפשריים מאתגרת הנחות ותיקות, ומסמנת מאמר זה כתרומה משמעותית לחקר התודעה.'
$d.Paragraphs(35).Range.Text = 'https://arxiv.org/abs/2409.14545'
